$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$a8 = $ws.Range("A8")
$a8.Characters().Text = $a8.Characters().Text.Replace("44", "45")

$c9 = $ws.Range("C9")
$c9.Characters().Text = $c9.Characters().Text.Replace("10/31/2022", "11/7/2022").Replace("11/6/2022", "11/13/2022")

# --- Donor cells (unchanged rows) used to copy exact cell style/type so no new styles are created ---
$donorText0  = $ws.Range("C14")  # style 14, shared text "0"
$donorTextNA = $ws.Range("E14")  # style 14, shared text "***.*"
$donorNum    = $ws.Range("I14")  # style 15, plain integer number format
$donorPct    = $ws.Range("K14")  # style 16, percent number format

# --- Cells changing data type (numeric <-> text placeholder) ---
$donorNum.Copy($ws.Range("C18"))
$ws.Range("C18").Value = 1
$donorNum.Copy($ws.Range("C22"))
$ws.Range("C22").Value = 2
$donorText0.Copy($ws.Range("D22"))
$donorTextNA.Copy($ws.Range("E22"))
$donorNum.Copy($ws.Range("C26"))
$ws.Range("C26").Value = 1
$donorText0.Copy($ws.Range("C27"))
$donorNum.Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$donorPct.Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$donorText0.Copy($ws.Range("D28"))
$donorTextNA.Copy($ws.Range("E28"))
$donorText0.Copy($ws.Range("F28"))
$donorText0.Copy($ws.Range("D29"))
$donorTextNA.Copy($ws.Range("E29"))
$donorText0.Copy($ws.Range("F29"))

# --- Plain value updates (style/type unchanged) ---
$ws.Range("M15").Value = 5.555555555555
$ws.Range("N15").Value = -34.482758620689
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 84.615384615384
$ws.Range("I16").Value = 201
$ws.Range("J16").Value = 166
$ws.Range("K16").Value = 21.084337349397
$ws.Range("L16").Value = 30.519480519480
$ws.Range("M16").Value = -13.733905579399
$ws.Range("N16").Value = -73.759791122715
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = -50
$ws.Range("I17").Value = 306
$ws.Range("J17").Value = 304
$ws.Range("K17").Value = 0.657894736842
$ws.Range("L17").Value = 23.387096774193
$ws.Range("M17").Value = 133.587786259542
$ws.Range("N17").Value = -16.621253405994
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 134
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = 44.086021505376
$ws.Range("L18").Value = -2.189781021897
$ws.Range("M18").Value = -47.450980392156
$ws.Range("N18").Value = -86.939571150097
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 36.363636363636
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = 3.773584905660
$ws.Range("I19").Value = 576
$ws.Range("J19").Value = 440
$ws.Range("K19").Value = 30.909090909090
$ws.Range("L19").Value = 22.033898305084
$ws.Range("M19").Value = 90.099009900990
$ws.Range("N19").Value = 10.982658959537
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -40
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = -8.333333333333
$ws.Range("I20").Value = 195
$ws.Range("J20").Value = 170
$ws.Range("K20").Value = 14.705882352941
$ws.Range("L20").Value = 20.370370370370
$ws.Range("M20").Value = -27.238805970149
$ws.Range("N20").Value = -93.155493155493
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -24.324324324324
$ws.Range("F21").Value = 130
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = -10.958904109589
$ws.Range("I21").Value = 1435
$ws.Range("J21").Value = 1189
$ws.Range("K21").Value = 20.689655172413
$ws.Range("L21").Value = 21.097046413502
$ws.Range("M21").Value = 18.399339933993
$ws.Range("N21").Value = -74.260089686098
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 24
$ws.Range("K22").Value = 71.428571428571
$ws.Range("L22").Value = 50
$ws.Range("M22").Value = 0
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 57.142857142857
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = -19.626168224299
$ws.Range("I24").Value = 1202
$ws.Range("J24").Value = 825
$ws.Range("K24").Value = 45.696969696969
$ws.Range("L24").Value = 63.093622795115
$ws.Range("M24").Value = 99.336650082918
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -9.090909090909
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 21.951219512195
$ws.Range("I25").Value = 467
$ws.Range("J25").Value = 412
$ws.Range("K25").Value = 13.349514563106
$ws.Range("L25").Value = 18.527918781725
$ws.Range("M25").Value = 5.656108597285
$ws.Range("F26").Value = 4
$ws.Range("I26").Value = 34
$ws.Range("K26").Value = 100
$ws.Range("L26").Value = 88.888888888888
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 400
$ws.Range("I27").Value = 52
$ws.Range("J27").Value = 35
$ws.Range("K27").Value = 48.571428571428
$ws.Range("L27").Value = 85.714285714285
$ws.Range("H28").Value = -100
$ws.Range("H29").Value = -100
